$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for years 2021 (col S) and 2022 (col T), rows 4-14.
# Column R already holds 2020 values with the correct formatting for each
# row; copy that formatting into S and T before writing the new values.
$values = @{
    4  = @(2021, 2022)
    5  = @(2.5, 2.6)
    6  = @(2.5, 1.8)
    7  = @(1.6, 2.6)
    8  = @(3.6, 1.9)
    9  = @(5.8, 3.9)
    10 = @(1.1000000000000001, 3.2)
    11 = @(1.1000000000000001, 3.3)
    12 = @(5.0999999999999996, 2.5)
    13 = @(2.2999999999999998, 1.9)
    14 = @(2.1, 2.5)
}

foreach ($row in 4..14) {
    $srcCell = $ws.Range("R$row")
    $dstRange = $ws.Range("S$row" + ":T$row")

    # Copy formatting (number format, fonts, borders -> style index) from R
    # into S:T for this row.
    $srcCell.Copy($dstRange)

    $pair = $values[$row]
    $ws.Range("S$row").Value = $pair[0]
    $ws.Range("T$row").Value = $pair[1]
}

# Update the saved selection to match the post-edit state.
[void]$ws.Range("V7").Select()
